# IR_export_translations.xlsx update:
# Insert two new translation rows ("IR_Screen_FLAG" / "IR_DataPrep_FLAG")
# into Sheet1 just above the existing "IR_Cat" row, shifting the rows
# below down by two, and update the view/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank rows at row 155 (pushes old rows 155-165 -> 157-167)
$ws.Rows("155:156").Insert()

# New row 155: IR_Screen_FLAG / DS
$ws.Cells.Item(155, 1).Value = "IR_Screen_FLAG"
$ws.Cells.Item(155, 2).Value = "DS"

# New row 156: IR_DataPrep_FLAG / DS
$ws.Cells.Item(156, 1).Value = "IR_DataPrep_FLAG"
$ws.Cells.Item(156, 2).Value = "DS"

# Update the view to match the committed state.
$ws.Application.ActiveWindow.ScrollRow = 149
$ws.Range("A156").Select()
